# Updated cryptos list on Sun May 19 20:54:34 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, plus
# a rank swap between VeChain (row 44) and dogwifhat (row 45).
#
# Column D holds prices as literal text (e.g. "3.075.89", "66.188.11")
# so the grouping dots survive round-tripping through Excel. Plain
# decimals in that column (e.g. "2.50", "0.999") would otherwise be
# auto-coerced to Double by the COM value setter and lose their exact
# text (2.50 -> 2.5), so every D-column write is prefixed with a leading
# apostrophe - the same trick as typing '2.50 directly into a cell -
# to force literal text storage. Column E's values already carry
# padding spaces ("  -1.23%  ") which keep them safely text without
# needing the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.188.11'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '''3.075.89'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''574.11'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').Value = '''170.25'
$ws.Range('E6').Value = '  -1.75%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''3.072.60'
$ws.Range('E8').Value = '  -1.60%  '
$ws.Range('E9').Value = '  -2.42%  '
$ws.Range('D10').Value = '''6.32'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('E12').Value = '  -2.71%  '
$ws.Range('D13').Value = '''0.0000239'
$ws.Range('E13').Value = '  -3.71%  '
$ws.Range('D14').Value = '''35.85'
$ws.Range('E14').Value = '  -3.73%  '
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '''3.588.57'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = '''66.156.12'
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('E18').Value = '  -2.96%  '
$ws.Range('D19').Value = '''3.073.17'
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('D20').Value = '''16.54'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').Value = '''484.21'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').Value = '''0.686'
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('D24').Value = '''82.39'
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('D25').Value = '''12.62'
$ws.Range('E25').Value = '  -5.28%  '
$ws.Range('E26').Value = '  -4.24%  '
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '''7.89'
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('E30').Value = '  -5.62%  '
$ws.Range('E31').Value = '  -3.29%  '
$ws.Range('D32').Value = '''27.72'
$ws.Range('E32').Value = '  -3.21%  '
$ws.Range('E33').Value = '  -4.47%  '
$ws.Range('D34').Value = '''0.0₃0898'
$ws.Range('E34').Value = '  -6.10%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D37').Value = '''47.13'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('E38').Value = '  -5.37%  '
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').Value = '''1.96'
$ws.Range('E40').Value = '  -5.39%  '
$ws.Range('D41').Value = '''0.299'
$ws.Range('E41').Value = '  -4.58%  '
$ws.Range('D42').Value = '''8.25'
$ws.Range('E42').Value = '  -4.99%  '
$ws.Range('D43').Value = '''2.778.33'
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '''2.50'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0343'
$ws.Range('E45').Value = '  -3.33%  '
$ws.Range('D46').Value = '''134.73'
$ws.Range('E46').Value = '  -0.74%  '
$ws.Range('D47').Value = '''363.68'
$ws.Range('E47').Value = '  -5.35%  '
$ws.Range('D49').Value = '''24.38'
$ws.Range('E49').Value = '  -2.44%  '
$ws.Range('D50').Value = '''2.15'
$ws.Range('E50').Value = '  -2.54%  '
$ws.Range('E51').Value = '  -2.56%  '
